$d = $word.ActiveDocument

$pairs = @(
  @{old="58+19="; new="34+27="},
  @{old="12+67="; new="37-36="},
  @{old="48-10="; new="0+41="},
  @{old="44-42="; new="13+5="},
  @{old="30-11="; new="35+29="},
  @{old="99-85="; new="18-15="},
  @{old="15+64="; new="15+6="},
  @{old="97-40="; new="75-12="},
  @{old="50+8="; new="29+55="},
  @{old="25-5="; new="11+83="},
  @{old="73-23="; new="57+17="},
  @{old="24+31="; new="91-87="},
  @{old="63-45="; new="86+5="},
  @{old="4+39="; new="38+40="},
  @{old="10-1="; new="45+33="},
  @{old="39+22="; new="36-28="},
  @{old="6+8="; new="33-29="},
  @{old="38+41="; new="41-34="},
  @{old="84-59="; new="2+92="},
  @{old="97-26="; new="3+0="},
  @{old="40+6="; new="98-36="},
  @{old="70+12="; new="93-41="},
  @{old="53+0="; new="48-18="},
  @{old="61-36="; new="89+8="},
  @{old="40+44="; new="89-16="},
  @{old="31-5="; new="10+25="},
  @{old="89-80="; new="16+51="},
  @{old="75-2="; new="60-44="},
  @{old="20+66="; new="37+36="},
  @{old="42-11="; new="55-38="},
  @{old="22-20="; new="20+57="},
  @{old="74-8="; new="65-43="},
  @{old="50+40="; new="95-26="},
  @{old="53+6="; new="12+38="},
  @{old="61-49="; new="77-19="},
  @{old="55-39="; new="3+43="},
  @{old="49+13="; new="20+46="},
  @{old="1+85="; new="44+48="},
  @{old="92-91="; new="24+51="},
  @{old="42-21="; new="15+47="},
  @{old="46-29="; new="5-5="},
  @{old="93-92="; new="42+2="},
  @{old="45-29="; new="24+17="},
  @{old="55+40="; new="16+0="},
  @{old="0+46="; new="95-58="},
  @{old="46+40="; new="34+45="},
  @{old="2+11="; new="41-30="},
  @{old="31-6="; new="96-59="},
  @{old="76-19="; new="31-2="},
  @{old="54+40="; new="89-24="},
  @{old="79-36="; new="34-33="},
  @{old="43-9="; new="20+31="},
  @{old="41-26="; new="90-24="},
  @{old="94-16="; new="63-24="},
  @{old="26+52="; new="72+3="},
  @{old="84-70="; new="70-9="},
  @{old="42-19="; new="56+13="},
  @{old="16+66="; new="34+59="},
  @{old="28+69="; new="74-65="},
  @{old="73-57="; new="0+31="},
  @{old="55+29="; new="18+37="},
  @{old="84-83="; new="68+27="},
  @{old="58-46="; new="9+37="},
  @{old="75+10="; new="14-13="},
  @{old="10+8="; new="40-24="},
  @{old="67-8="; new="69-5="},
  @{old="7+66="; new="71-10="},
  @{old="16+73="; new="38+28="},
  @{old="31+6="; new="18+64="},
  @{old="79+2="; new="62+23="},
  @{old="27-3="; new="91-82="},
  @{old="18+50="; new="59+38="},
  @{old="69-32="; new="57+37="},
  @{old="22+5="; new="65-60="},
  @{old="37+58="; new="10+56="},
  @{old="5+74="; new="45-9="},
  @{old="67+11="; new="35+13="},
  @{old="1+37="; new="63-25="},
  @{old="6+57="; new="88-44="},
  @{old="84-58="; new="66-7="},
  @{old="40+19="; new="9+87="},
  @{old="25+25="; new="5+20="},
  @{old="88-18="; new="11+34="},
  @{old="73-28="; new="27-10="},
  @{old="78-49="; new="35-13="},
  @{old="2+97="; new="67-29="},
  @{old="14+74="; new="38+54="},
  @{old="24+48="; new="70-26="},
  @{old="77-14="; new="31-4="},
  @{old="13+65="; new="7+59="},
  @{old="73-33="; new="32-31="},
  @{old="18-14="; new="50+18="},
  @{old="84-2="; new="78+12="},
  @{old="50+21="; new="2+16="},
  @{old="71-35="; new="82-72="},
  @{old="53-28="; new="35+35="},
  @{old="46+10="; new="76-65="},
  @{old="41-23="; new="86+8="},
  @{old="87-8="; new="15+0="},
  @{old="64+2="; new="97-45="}
)

foreach ($pair in $pairs) {
  $range = $d.Content
  $range.Find.Execute($pair.old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}